# Prakruti timesheet wk 14 - Signed Off time sheets
# Fill in the Supervisor Name, and the Supervisor Signature / Date fields
# (the employee side was already signed off; this records the supervisor's
# sign-off).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name (merged G6:I6)
$ws.Range("G6").Value = "Ankita Gangotra"

# Supervisor Signature initials (merged A27:C27)
$ws.Range("A27").Value = "A.G"

# Supervisor Signature date (merged D27:E27) - 10 June 2014
$ws.Range("D27").Value = 41800
$ws.Range("D27").NumberFormat = "m/d/yy"

# Leave the selection on the date field just filled in, matching the
# final UI state after the sign-off edit.
$ws.Range("D27:E27").Select()
